# Weekly update: insert a new data row (new week) at row 22, pushing
# the existing rows 22-40 down to 23-41 (dimension grows from A1:R40 to A1:R41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22, shifting rows 22..40 down to 23..41.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with this week's record.
$ws.Cells.Item(22, 1).Value = 3
$ws.Cells.Item(22, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44873
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = 300000000
$ws.Cells.Item(22, 7).Value = "Espárragos"
$ws.Cells.Item(22, 8).Value = "Verde"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 11).Value = 1400
$ws.Cells.Item(22, 12).Value = 1450
$ws.Cells.Item(22, 13).Value = 1423
$ws.Cells.Item(22, 14).Value = "`$/kilo"
$ws.Cells.Item(22, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(22, 16).Value = 1423
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
